$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("quality_comparison")
$ws2 = $wb.Worksheets.Item("computational_comparison")

# --- New partial-border cell styles ---
# "top+bottom" template, built once on sheet1!C1
$ws1.Range("C1").ClearFormats()
$ws1.Range("C1").Borders.Item(8).LineStyle = 1    # xlEdgeTop
$ws1.Range("C1").Borders.Item(9).LineStyle = 1    # xlEdgeBottom

# "top+right+bottom" template, built once on sheet1!D1
$ws1.Range("D1").ClearFormats()
$ws1.Range("D1").Borders.Item(8).LineStyle = 1    # xlEdgeTop
$ws1.Range("D1").Borders.Item(10).LineStyle = 1   # xlEdgeRight
$ws1.Range("D1").Borders.Item(9).LineStyle = 1    # xlEdgeBottom

# Propagate the two new styles onto sheet2's matching header cells
# (copy/paste-format reuses the same style index instead of minting new ones)
$ws1.Range("C1").Copy()
$ws2.Range("C1").PasteSpecial(-4122)   # xlPasteFormats
$ws2.Range("F1").PasteSpecial(-4122)

$ws1.Range("D1").Copy()
$ws2.Range("D1").PasteSpecial(-4122)
$ws2.Range("G1").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Anonymize "fedcore" -> "approach" ---
$ws1.Range("C2").Value = "approach"
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# --- Drop the stray empty inline-string cell G5 on sheet2 ---
$ws2.Range("G5").ClearContents()
